# The "특102" sheet used to hold a small purchase/expense table
# (삼겹살/식탁보/치즈김밥/김밥 rows). That data actually belonged on the
# "빈소2" sheet, so move it there - trimming it down to just the two
# rows that are still relevant (치즈김밥 x24, 식탁보 x1) - and clear the
# old table out of "특102" so it goes back to being a blank sheet.

$wb = $excel.ActiveWorkbook

$dest = $wb.Worksheets.Item("빈소2")
$dest.Range("A1").Value = "치즈김밥"
$dest.Range("B1").Value = 2500
$dest.Range("C1").Value = 24
$dest.Range("D1").Value = 60000
$dest.Range("E1").Value = " "

$dest.Range("A4").Value = "식탁보"
$dest.Range("B4").Value = 6000
$dest.Range("C4").Value = 1
$dest.Range("D4").Value = 6000
$dest.Range("E4").Value = " "

$src = $wb.Worksheets.Item("특102")
$src.Range("A1:E4").ClearContents()
